# Update "want to go" (想去人数) counts on the 展览, 演出 and 全部类型 sheets.
# (本地生活 has no changes in this edit.)

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 167
$ws.Range("F3").Value  = 1403
$ws.Range("F7").Value  = 522
$ws.Range("F8").Value  = 853
$ws.Range("F9").Value  = 571
$ws.Range("F10").Value = 775
$ws.Range("F11").Value = 349
$ws.Range("F14").Value = 1099
$ws.Range("F15").Value = 537
$ws.Range("F16").Value = 314
$ws.Range("F19").Value = 270
$ws.Range("F21").Value = 62
$ws.Range("F23").Value = 493
$ws.Range("F25").Value = 409

# --- 演出 sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 52
$ws.Range("F5").Value  = 301
$ws.Range("F10").Value = 176
$ws.Range("F11").Value = 39

# --- 全部类型 sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 167
$ws.Range("F4").Value  = 1403
$ws.Range("F10").Value = 52
$ws.Range("F11").Value = 301
$ws.Range("F12").Value = 522
$ws.Range("F13").Value = 853
$ws.Range("F14").Value = 571
$ws.Range("F15").Value = 775
$ws.Range("F16").Value = 349
$ws.Range("F19").Value = 1099
$ws.Range("F20").Value = 537
$ws.Range("F23").Value = 314
$ws.Range("F28").Value = 270
$ws.Range("F30").Value = 62
$ws.Range("F31").Value = 176
$ws.Range("F32").Value = 39
$ws.Range("F36").Value = 493
$ws.Range("F38").Value = 409
